$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# Widen column B (target OOXML width="24"; the engine's char->XML-unit
# conversion needs a slightly larger ColumnWidth input to land exactly on 24)
$ws.Columns.Item(2).ColumnWidth = 23.25

# Enter the new row of data (row 8), in the same order the original author
# typed it in (A, then C/Description, then B/Value) so the shared-string
# table gets the same append order as the canonical file.
$ws.Range("A8").Value = "Header Format"
$ws.Range("C8").Value = "Must be TIME::VALUE_TYPE"
$ws.Range("B8").Value = "TIME::VALUE_TYPE"

# Match the formatting used by the other rows: column A / C use the style
# from the row above, column B (Value) uses the same style as B2 (the only
# other populated Value cell).
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# New row keeps the sheet's default row height, but explicitly so.
$ws.Rows.Item(8).RowHeight = 18

# Selection moves to the newly entered cell.
$ws.Range("C8").Select() | Out-Null

# Page setup was touched (paper size + portrait orientation).
$ws.PageSetup.PaperSize = 10
$ws.PageSetup.Orientation = 1
